$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 4372
$ws1.Range("F10").Value = 786
$ws1.Range("F24").Value = 2763
$ws1.Range("F26").Value = 293

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 4372
$ws4.Range("F10").Value = 786
$ws4.Range("F25").Value = 2763
$ws4.Range("F27").Value = 293
